# Auto-generated Excel COM-interop script to apply numeric corrections
# to the Phantom_Profits workbook tables across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 817.58826
$ws.Range("I9").Value = 842.7857
$ws.Range("K9").Value = 842.7857
$ws.Range("M9").Value = -673.7857
$ws.Range("H28").Value = 1427.619
$ws.Range("I28").Value = 1117.4375
$ws.Range("J28").Value = 2420.2
$ws.Range("K28").Value = 1117.4375
$ws.Range("L28").Value = 2420.2
$ws.Range("M28").Value = -632.4375
$ws.Range("N28").Value = -3390.2
$ws.Range("H32").Value = 7330
$ws.Range("I32").Value = 3001
$ws.Range("J32").Value = 9494.5
$ws.Range("K32").Value = 3001
$ws.Range("L32").Value = 9494.5
$ws.Range("M32").Value = -2675
$ws.Range("N32").Value = -10146.5
$ws.Range("H92").Value = 36145.57
$ws.Range("I92").Value = 59317.06
$ws.Range("J92").Value = 335.0909
$ws.Range("K92").Value = 59317.06
$ws.Range("L92").Value = 335.0909
$ws.Range("M92").Value = -58069.06
$ws.Range("N92").Value = -2831.0909
$ws.Range("H121").Value = 2083.3333
$ws.Range("J121").Value = 2083.3333
$ws.Range("L121").Value = 6249.999899999999
$ws.Range("N121").Value = -9743.999899999999
$ws.Range("H132").Value = 4783.5
$ws.Range("I132").Value = 4783.5
$ws.Range("K132").Value = 14350.5
$ws.Range("M132").Value = -11820.5
$ws.Range("H135").Value = 685.6111
$ws.Range("I135").Value = 637.625
$ws.Range("K135").Value = 5738.625
$ws.Range("M135").Value = -3203.625
$ws.Range("H137").Value = 3000.7646
$ws.Range("I137").Value = 1251.4
$ws.Range("K137").Value = 3754.2
$ws.Range("M137").Value = -1204.2
$ws.Range("H138").Value = 1914
$ws.Range("I138").Value = 1733.409
$ws.Range("K138").Value = 5200.227000000001
$ws.Range("M138").Value = -60.22700000000077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3436.6
$ws.Range("I32").Value = 2727.3794
$ws.Range("K32").Value = 2727.3794
$ws.Range("M32").Value = -2440.3794
$ws.Range("H44").Value = 65000
$ws.Range("J44").Value = 65000
$ws.Range("L44").Value = 65000
$ws.Range("N44").Value = -65976
$ws.Range("H61").Value = 5924.5557
$ws.Range("I61").Value = 5245.5
$ws.Range("J61").Value = 8301.25
$ws.Range("K61").Value = 5245.5
$ws.Range("L61").Value = 8301.25
$ws.Range("M61").Value = -5033.5
$ws.Range("N61").Value = -8725.25
$ws.Range("H88").Value = 2755.2307
$ws.Range("I88").Value = 2037.5
$ws.Range("K88").Value = 2037.5
$ws.Range("M88").Value = -1631.5
$ws.Range("H91").Value = 2755.2307
$ws.Range("I91").Value = 2037.5
$ws.Range("K91").Value = 2037.5
$ws.Range("M91").Value = -633.5
$ws.Range("H120").Value = 112500
$ws.Range("J120").Value = 112500
$ws.Range("L120").Value = 112500
$ws.Range("N120").Value = -122176
$ws.Range("H122").Value = 1761.1666
$ws.Range("I122").Value = 1687.4
$ws.Range("J122").Value = 2130
$ws.Range("K122").Value = 5062.200000000001
$ws.Range("L122").Value = 6390
$ws.Range("M122").Value = -2612.200000000001
$ws.Range("N122").Value = -11290
$ws.Range("H133").Value = 59397.5
$ws.Range("J133").Value = 59397.5
$ws.Range("L133").Value = 59397.5
$ws.Range("N133").Value = -64457.5
$ws.Range("H136").Value = 5924.5557
$ws.Range("I136").Value = 5245.5
$ws.Range("J136").Value = 8301.25
$ws.Range("K136").Value = 15736.5
$ws.Range("L136").Value = 24903.75
$ws.Range("M136").Value = -13186.5
$ws.Range("N136").Value = -30003.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3071.7144
$ws.Range("J20").Value = 2999.5
$ws.Range("L20").Value = 2999.5
$ws.Range("N20").Value = -3493.5
$ws.Range("H86").Value = 11799.308
$ws.Range("I86").Value = 4577.4
$ws.Range("K86").Value = 4577.4
$ws.Range("M86").Value = -3454.4
$ws.Range("H89").Value = 11799.308
$ws.Range("I89").Value = 4577.4
$ws.Range("K89").Value = 22887
$ws.Range("M89").Value = -17271
$ws.Range("H107").Value = 2600.25
$ws.Range("I107").Value = 1955.5
$ws.Range("J107").Value = 3245
$ws.Range("K107").Value = 1955.5
$ws.Range("L107").Value = 3245
$ws.Range("M107").Value = -35.5
$ws.Range("N107").Value = -7085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4637
$ws.Range("I31").Value = 2862.889
$ws.Range("J31").Value = 9959.333000000001
$ws.Range("K31").Value = 2862.889
$ws.Range("L31").Value = 9959.333000000001
$ws.Range("M31").Value = -2567.889
$ws.Range("N31").Value = -10549.333
$ws.Range("H34").Value = 4637
$ws.Range("I34").Value = 2862.889
$ws.Range("J34").Value = 9959.333000000001
$ws.Range("K34").Value = 2862.889
$ws.Range("L34").Value = 9959.333000000001
$ws.Range("M34").Value = -2660.889
$ws.Range("N34").Value = -10363.333
$ws.Range("H86").Value = 3083.1667
$ws.Range("I86").Value = 3104.889
$ws.Range("K86").Value = 3104.889
$ws.Range("M86").Value = -1981.889
$ws.Range("H89").Value = 3083.1667
$ws.Range("I89").Value = 3104.889
$ws.Range("K89").Value = 15524.445
$ws.Range("M89").Value = -9908.445
$ws.Range("H132").Value = 1703.0588
$ws.Range("J132").Value = 1857
$ws.Range("L132").Value = 5571
$ws.Range("N132").Value = -10631
$ws.Range("H141").Value = 71300.60000000001
$ws.Range("J141").Value = 95402.336
$ws.Range("L141").Value = 95402.336
$ws.Range("N141").Value = -105762.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H121").Value = 911765.8
$ws.Range("I121").Value = 1121.5
$ws.Range("J121").Value = 2004539
$ws.Range("K121").Value = 3364.5
$ws.Range("L121").Value = 6013617
$ws.Range("M121").Value = -2054.5
$ws.Range("N121").Value = -6016237
$ws.Range("H132").Value = 1531
$ws.Range("J132").Value = 1733.3334
$ws.Range("L132").Value = 15600.0006
$ws.Range("N132").Value = -20660.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6490.1665
$ws.Range("J107").Value = 9665.666999999999
$ws.Range("L107").Value = 9665.666999999999
$ws.Range("N107").Value = -13505.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 760.8
$ws.Range("I22").Value = 468.33334
$ws.Range("J22").Value = 1199.5
$ws.Range("K22").Value = 468.33334
$ws.Range("L22").Value = 1199.5
$ws.Range("M22").Value = -173.33334
$ws.Range("N22").Value = -1789.5
$ws.Range("H27").Value = 760.8
$ws.Range("I27").Value = 468.33334
$ws.Range("J27").Value = 1199.5
$ws.Range("K27").Value = 468.33334
$ws.Range("L27").Value = 1199.5
$ws.Range("M27").Value = -361.33334
$ws.Range("N27").Value = -1413.5
$ws.Range("H55").Value = 452.57144
$ws.Range("I55").Value = 339.66666
$ws.Range("J55").Value = 537.25
$ws.Range("K55").Value = 339.66666
$ws.Range("L55").Value = 537.25
$ws.Range("M55").Value = -166.66666
$ws.Range("N55").Value = -883.25
$ws.Range("H100").Value = 2034
$ws.Range("I100").Value = 2001
$ws.Range("K100").Value = 2001
$ws.Range("M100").Value = -1460
$ws.Range("H132").Value = 2324.0454
$ws.Range("I132").Value = 2275.3157
$ws.Range("K132").Value = 6825.9471
$ws.Range("M132").Value = -4295.9471

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 487.5263
$ws.Range("I113").Value = 365.2857
$ws.Range("J113").Value = 829.8
$ws.Range("K113").Value = 1095.8571
$ws.Range("L113").Value = 2489.4
$ws.Range("M113").Value = 1074.1429
$ws.Range("N113").Value = -6829.4
$ws.Range("H122").Value = 5606.5557
$ws.Range("I122").Value = 5137.8667
$ws.Range("J122").Value = 7950
$ws.Range("K122").Value = 15413.6001
$ws.Range("L122").Value = 23850
$ws.Range("M122").Value = -12963.6001
$ws.Range("N122").Value = -28750
